$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $s = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $s
}

$sub3 = [char]8323
$d40val = [string]::Concat('0.0', $sub3, '0737')

Set-TextValue "D2" '69.169.69'
Set-TextValue "E2" '  +1.22%  '
Set-TextValue "D3" '3.405.85'
Set-TextValue "E3" '  +1.59%  '
Set-TextValue "E4" '  -0.06%  '
Set-TextValue "D5" '581.94'
Set-TextValue "E5" '  -0.51%  '
Set-TextValue "D6" '178.44'
Set-TextValue "E6" '  +0.69%  '
Set-TextValue "E7" '  +0.07%  '
Set-TextValue "E8" '  +0.48%  '
Set-TextValue "E9" '  +7.70%  '
Set-TextValue "D10" '0.586'
Set-TextValue "E10" '  +0.69%  '
Set-TextValue "D11" '48.39'
Set-TextValue "E11" '  +0.94%  '
Set-TextValue "E12" '  +3.00%  '
Set-TextValue "D13" '679.16'
Set-TextValue "E13" '  -1.47%  '
Set-TextValue "D14" '3.950.41'
Set-TextValue "E14" '  +1.15%  '
Set-TextValue "D15" '8.61'
Set-TextValue "E15" '  +1.98%  '
Set-TextValue "D16" '69.244.62'
Set-TextValue "E16" '  +1.28%  '
Set-TextValue "D17" '3.400.92'
Set-TextValue "E17" '  +1.21%  '
Set-TextValue "E18" '  +0.47%  '
Set-TextValue "D19" '17.73'
Set-TextValue "E19" '  +1.44%  '
Set-TextValue "D20" '11.29'
Set-TextValue "E20" '  +0.89%  '
Set-TextValue "D21" '0.911'
Set-TextValue "E21" '  +1.87%  '
Set-TextValue "D22" '5.39'
Set-TextValue "E22" '  -1.37%  '
Set-TextValue "D23" '17.01'
Set-TextValue "E23" '  +0.46%  '
Set-TextValue "D24" '100.65'
Set-TextValue "E24" '  +0.56%  '
Set-TextValue "E25" '  -0.34%  '
Set-TextValue "E27" '  +1.82%  '
Set-TextValue "D28" '33.46'
Set-TextValue "E28" '  +1.50%  '
Set-TextValue "E29" '  +2.58%  '
Set-TextValue "D30" '6.86'
Set-TextValue "E30" '  -1.21%  '
Set-TextValue "E31" '  +10.48%  '
Set-TextValue "D32" '555.63'
Set-TextValue "E32" '  +0.47%  '
Set-TextValue "E33" '  -0.74%  '
Set-TextValue "E34" '  -0.39%  '
Set-TextValue "D35" '58.03'
Set-TextValue "E35" '  +0.01%  '
Set-TextValue "E36" '  +0.09%  '
Set-TextValue "D37" '3.608.36'
Set-TextValue "E37" '  -2.98%  '
Set-TextValue "E38" '  +0.32%  '
Set-TextValue "D39" '34.92'
Set-TextValue "E39" '  +0.71%  '
Set-TextValue "D40" $d40val
Set-TextValue "E40" '  +9.85%  '
Set-TextValue "D41" '3.28'
Set-TextValue "E41" '  +3.32%  '
Set-TextValue "E42" '  +2.90%  '
Set-TextValue "E43" '  +2.62%  '
Set-TextValue "D44" '0.0424'
Set-TextValue "E44" '  +3.23%  '
Set-TextValue "D45" '0.334'
Set-TextValue "E46" '  +0.80%  '
Set-TextValue "E47" '  +0.28%  '
Set-TextValue "E48" '  +3.64%  '
Set-TextValue "E49" '  -0.07%  '
Set-TextValue "D50" '131.04'
Set-TextValue "E50" '  -0.59%  '
Set-TextValue "D51" '2.67'
Set-TextValue "E51" '  +3.19%  '
